# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market-data refresh to the Leve-profit sheets.
# For each (sheet, row) pair below, columns H-N are overwritten with freshly
# fetched market-price derived values; a couple of rows also gain/lose a cell
# (M/N) entirely, matching the upstream diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4683.3335
$ws.Cells.Item(64, 9).Value = 6044.4443
$ws.Cells.Item(64, 10).Value = 3322.2222
$ws.Cells.Item(64, 11).Value = 6044.4443
$ws.Cells.Item(64, 12).Value = 3322.2222
$ws.Cells.Item(64, 13).Value = -5796.4443
$ws.Cells.Item(64, 14).Value = -3818.2222

$ws.Cells.Item(67, 8).Value = 4683.3335
$ws.Cells.Item(67, 9).Value = 6044.4443
$ws.Cells.Item(67, 10).Value = 3322.2222
$ws.Cells.Item(67, 11).Value = 6044.4443
$ws.Cells.Item(67, 12).Value = 3322.2222
$ws.Cells.Item(67, 13).Value = -5186.4443
$ws.Cells.Item(67, 14).Value = -5038.2222

$ws.Cells.Item(76, 8).Value = 6671359
$ws.Cells.Item(76, 9).Value = 11908748
$ws.Cells.Item(76, 10).Value = 5590.909
$ws.Cells.Item(76, 11).Value = 11908748
$ws.Cells.Item(76, 12).Value = 5590.909
$ws.Cells.Item(76, 13).Value = -11908433
$ws.Cells.Item(76, 14).Value = -6220.909

$ws.Cells.Item(79, 8).Value = 6671359
$ws.Cells.Item(79, 9).Value = 11908748
$ws.Cells.Item(79, 10).Value = 5590.909
$ws.Cells.Item(79, 11).Value = 11908748
$ws.Cells.Item(79, 12).Value = 5590.909
$ws.Cells.Item(79, 13).Value = -11907656
$ws.Cells.Item(79, 14).Value = -7774.909

$ws.Cells.Item(129, 8).Value = 1157.2142
$ws.Cells.Item(129, 9).Value = 802.4545000000001
$ws.Cells.Item(129, 10).Value = 1386.7646
$ws.Cells.Item(129, 11).Value = 2407.3635
$ws.Cells.Item(129, 12).Value = 4160.293799999999
$ws.Cells.Item(129, 13).Value = 2592.6365
$ws.Cells.Item(129, 14).Value = -14160.2938

$ws.Cells.Item(138, 8).Value = 3113.2195
$ws.Cells.Item(138, 9).Value = 1504.7097
$ws.Cells.Item(138, 10).Value = 4090.9412
$ws.Cells.Item(138, 11).Value = 4514.1291
$ws.Cells.Item(138, 12).Value = 12272.8236
$ws.Cells.Item(138, 13).Value = 625.8708999999999
$ws.Cells.Item(138, 14).Value = -22552.8236

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9825.389999999999
$ws.Cells.Item(32, 9).Value = 7670.811
$ws.Cells.Item(32, 10).Value = 15957.654
$ws.Cells.Item(32, 11).Value = 7670.811
$ws.Cells.Item(32, 12).Value = 15957.654
$ws.Cells.Item(32, 13).Value = -7383.811
$ws.Cells.Item(32, 14).Value = -16531.654

$ws.Cells.Item(63, 8).Value = 166668530
$ws.Cells.Item(63, 9).Value = 200001920
$ws.Cells.Item(63, 10).Value = 1600
$ws.Cells.Item(63, 11).Value = 200001920
$ws.Cells.Item(63, 12).Value = 1600
$ws.Cells.Item(63, 13).Value = -200001234
$ws.Cells.Item(63, 14).Value = -2972

$ws.Cells.Item(66, 8).Value = 166668530
$ws.Cells.Item(66, 9).Value = 200001920
$ws.Cells.Item(66, 10).Value = 1600
$ws.Cells.Item(66, 11).Value = 1000009600
$ws.Cells.Item(66, 12).Value = 8000
$ws.Cells.Item(66, 13).Value = -1000006168
$ws.Cells.Item(66, 14).Value = -14864

$ws.Cells.Item(74, 8).Value = 17858868
$ws.Cells.Item(74, 9).Value = 1431.1538
$ws.Cells.Item(74, 10).Value = 33335314
$ws.Cells.Item(74, 11).Value = 1431.1538
$ws.Cells.Item(74, 12).Value = 33335314
$ws.Cells.Item(74, 13).Value = -557.1538
$ws.Cells.Item(74, 14).Value = -33337062

$ws.Cells.Item(77, 8).Value = 17858868
$ws.Cells.Item(77, 9).Value = 1431.1538
$ws.Cells.Item(77, 10).Value = 33335314
$ws.Cells.Item(77, 11).Value = 7155.769
$ws.Cells.Item(77, 12).Value = 166676570
$ws.Cells.Item(77, 13).Value = -2787.769
$ws.Cells.Item(77, 14).Value = -166685306

$ws.Cells.Item(132, 8).Value = 5176.8647
$ws.Cells.Item(132, 9).Value = 2955.8635
$ws.Cells.Item(132, 10).Value = 8434.333000000001
$ws.Cells.Item(132, 11).Value = 8867.5905
$ws.Cells.Item(132, 12).Value = 25302.999
$ws.Cells.Item(132, 13).Value = -6337.5905
$ws.Cells.Item(132, 14).Value = -30362.999

$ws.Cells.Item(135, 8).Value = 47838.875
$ws.Cells.Item(135, 10).Value = 47838.875
$ws.Cells.Item(135, 12).Value = 47838.875
$ws.Cells.Item(135, 14).Value = -57978.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 49900
$ws.Cells.Item(81, 10).Value = 49900
$ws.Cells.Item(81, 12).Value = 49900
$ws.Cells.Item(81, 14).Value = -52022

$ws.Cells.Item(84, 8).Value = 49900
$ws.Cells.Item(84, 10).Value = 49900
$ws.Cells.Item(84, 12).Value = 149700
$ws.Cells.Item(84, 14).Value = -160308

$ws.Cells.Item(105, 8).Value = 3006.8572
$ws.Cells.Item(105, 9).Value = 1861.7
$ws.Cells.Item(105, 11).Value = 1861.7
$ws.Cells.Item(105, 13).Value = -114.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 10422314
$ws.Cells.Item(31, 9).Value = 1928.125
$ws.Cells.Item(31, 10).Value = 15632506
$ws.Cells.Item(31, 11).Value = 1928.125
$ws.Cells.Item(31, 12).Value = 15632506
$ws.Cells.Item(31, 13).Value = -1633.125
$ws.Cells.Item(31, 14).Value = -15633096

$ws.Cells.Item(34, 8).Value = 10422314
$ws.Cells.Item(34, 9).Value = 1928.125
$ws.Cells.Item(34, 10).Value = 15632506
$ws.Cells.Item(34, 11).Value = 1928.125
$ws.Cells.Item(34, 12).Value = 15632506
$ws.Cells.Item(34, 13).Value = -1726.125
$ws.Cells.Item(34, 14).Value = -15632910

$ws.Cells.Item(62, 8).Value = 17500
$ws.Cells.Item(62, 9).Value = 17375
$ws.Cells.Item(62, 11).Value = 17375
$ws.Cells.Item(62, 13).Value = -16751

$ws.Cells.Item(65, 8).Value = 17500
$ws.Cells.Item(65, 9).Value = 17375
$ws.Cells.Item(65, 11).Value = 86875
$ws.Cells.Item(65, 13).Value = -83755

$ws.Cells.Item(132, 8).Value = 6898359
$ws.Cells.Item(132, 9).Value = 11112698
$ws.Cells.Item(132, 11).Value = 33338094
$ws.Cells.Item(132, 13).Value = -33335564

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(125, 8).Value = 4686.1333
$ws.Cells.Item(125, 9).Value = 3171.8333
$ws.Cells.Item(125, 10).Value = 5695.6665
$ws.Cells.Item(125, 11).Value = 9515.499899999999
$ws.Cells.Item(125, 12).Value = 17086.9995
$ws.Cells.Item(125, 13).Value = -4595.499899999999
$ws.Cells.Item(125, 14).Value = -26926.9995

$ws.Cells.Item(131, 8).Value = 2174921
$ws.Cells.Item(131, 10).Value = 1144.2188
$ws.Cells.Item(131, 12).Value = 3432.6564
$ws.Cells.Item(131, 14).Value = -13512.6564

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4795.56
$ws.Cells.Item(70, 9).Value = 4527.7856
$ws.Cells.Item(70, 11).Value = 4527.7856
$ws.Cells.Item(70, 13).Value = -4257.7856

$ws.Cells.Item(73, 8).Value = 4795.56
$ws.Cells.Item(73, 9).Value = 4527.7856
$ws.Cells.Item(73, 11).Value = 4527.7856
$ws.Cells.Item(73, 13).Value = -3591.7856

$ws.Cells.Item(80, 8).Value = 13279
$ws.Cells.Item(80, 9).Value = 14638.125
$ws.Cells.Item(80, 10).Value = 2406
$ws.Cells.Item(80, 11).Value = 14638.125
$ws.Cells.Item(80, 12).Value = 2406
$ws.Cells.Item(80, 13).Value = -13640.125
$ws.Cells.Item(80, 14).Value = -4402

$ws.Cells.Item(83, 8).Value = 13279
$ws.Cells.Item(83, 9).Value = 14638.125
$ws.Cells.Item(83, 10).Value = 2406
$ws.Cells.Item(83, 11).Value = 73190.625
$ws.Cells.Item(83, 12).Value = 12030
$ws.Cells.Item(83, 13).Value = -68198.625
$ws.Cells.Item(83, 14).Value = -22014

$ws.Cells.Item(126, 8).Value = 12643.056
$ws.Cells.Item(126, 9).Value = 17382.5
$ws.Cells.Item(126, 10).Value = 3164.1667
$ws.Cells.Item(126, 11).Value = 52147.5
$ws.Cells.Item(126, 12).Value = 9492.500100000001
$ws.Cells.Item(126, 13).Value = -49677.5
$ws.Cells.Item(126, 14).Value = -14432.5001

$ws.Cells.Item(132, 8).Value = 5558969
$ws.Cells.Item(132, 9).Value = 10417679
$ws.Cells.Item(132, 11).Value = 31253037
$ws.Cells.Item(132, 13).Value = -31250507

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(128, 8).Value = 30390
$ws.Cells.Item(128, 9).Value = 30390
$ws.Cells.Item(128, 11).Value = 30390
$ws.Cells.Item(128, 13).Value = -25410

$ws.Cells.Item(129, 8).Value = 15000
$ws.Cells.Item(129, 9).Value = 15000
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 11).Value = 15000
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(129, 13).Value = -10000
$ws.Cells.Item(129, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 3286.6904
$ws.Cells.Item(132, 9).Value = 3061.3428
$ws.Cells.Item(132, 11).Value = 9184.028399999999
$ws.Cells.Item(132, 13).Value = -6654.028399999999

$ws.Cells.Item(133, 8).Value = 79441.664
$ws.Cells.Item(133, 10).Value = 79441.664
$ws.Cells.Item(133, 12).Value = 79441.664
$ws.Cells.Item(133, 14).Value = -84501.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1794.45
$ws.Cells.Item(132, 9).Value = 822.1111
$ws.Cells.Item(132, 10).Value = 2590
$ws.Cells.Item(132, 11).Value = 2466.3333
$ws.Cells.Item(132, 12).Value = 7770
$ws.Cells.Item(132, 13).Value = 63.66670000000022
$ws.Cells.Item(132, 14).Value = -12830

$ws.Cells.Item(136, 8).Value = 6670884.5
$ws.Cells.Item(136, 9).Value = 6271.2
$ws.Cells.Item(136, 10).Value = 11113960
$ws.Cells.Item(136, 11).Value = 18813.6
$ws.Cells.Item(136, 12).Value = 33341880
$ws.Cells.Item(136, 13).Value = -16263.6
$ws.Cells.Item(136, 14).Value = -33346980
